$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EffectiveDate (F) and PreviousExpDate (I) values from 11012023 to 12012023
# for every data row (2-9)
foreach ($r in 2..9) {
    $ws.Range("F$r").Value = "12012023"
    $ws.Range("I$r").Value = "12012023"
}

# Update Quality (P) from "Standard" to "Economy" for the "Teddy/George" rows (3,5,7,9)
foreach ($r in 3,5,7,9) {
    $ws.Range("P$r").Value = "Economy"
}

# Row 7 now has a custom height
$ws.Rows("7").RowHeight = 23.4

# Remove the last data row (row 9, Teddy/George) and the two trailing blank rows (11 & 12),
# shifting everything up so the sheet ends at row 9 (used to be row 10, a blank row).
$ws.Rows("9").Delete()
$ws.Rows("10:11").Delete()

# Update the active selection
$ws.Range("D11").Select() | Out-Null
